$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain numeric-looking values as literal
# text (e.g. "29.374.39", "6.390", "0.4340"). Force text format on each
# touched cell before assigning so Excel does not reinterpret/normalize
# them as numbers (which would drop trailing zeros or merge the dotted
# thousands groups).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.374.39'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.879.87'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7184'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.09'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07921'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3144'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08145'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.888.09'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '94.69'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.233'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7095'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.390'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008408'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.376.49'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '250.17'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.30'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.126.72'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1587'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.087'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.47'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.86'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.506'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.412'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.282'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.223'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05318'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.937'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7544'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.703'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01891'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.273.98'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.763'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.464'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '113.13'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '74.44'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9054'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.003'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.023.12'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.799'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.486'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4340'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.09%  '
